$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix grades for rows 29 and 30 (B29, B30) -- the dependent SUM formulas
# (B35 and B36) will recalculate automatically.
$ws.Range("B29").Value = 2
$ws.Range("B30").Value = 5

# Move the view/selection to match the saved sheet view state.
[void]$ws.Range("E29").Select()
